# Regenerate Report for Handoff:
# the e2e test markdown file's GUID-based name changed from
#   f884e64a-da2d-4733-b9fc-5711cae29439
# to
#   36db075a-336b-4ecf-ad40-2d245c6e7a1c
# and the handoff/handback xlf files picked up a new content hash
# (6e5e33d0f0f70fc3646c13ab3a5666357c283877 -> 84f50f49420392f2643eb20fea570acd54a163d3),
# plus refreshed "Ready for handoff" timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "f884e64a-da2d-4733-b9fc-5711cae29439"
$newGuid = "36db075a-336b-4ecf-ad40-2d245c6e7a1c"
$oldHash = "6e5e33d0f0f70fc3646c13ab3a5666357c283877"
$newHash = "84f50f49420392f2643eb20fea570acd54a163d3"

$newMdName = "$newGuid.md"
$newMdPath = "e2e\$newGuid.md"
$newHandoffDate = "2016-08-24 02:57:08"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"
$newZhHandoffDate = "2016-08-24 02:56:59"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# The external hyperlink target (rId) stays pointed at the original
# github blob URL - only the cell text / hyperlink display text changes.
$oldHyperlinkAddress = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/f5692af69490aecc47d50326f44da14ec7fa5539/e2e/$oldGuid.md"

# ---------------------------------------------------------------------
# Sheet "Overview": A2 file name, B2 path+name (hyperlink), G2 handoff date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = $newMdName
$wsOverview.Range("G2").Value = $newHandoffDate

$rB2 = $wsOverview.Range("B2")
$rB2.Hyperlinks.Delete()
$wsOverview.Hyperlinks.Add($rB2, $oldHyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newMdPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn": A2 file name (hyperlink), G2 latest handoff xlf file,
# H2 latest handoff datetime
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("G2").Value = $newZhXlf
$wsZh.Range("H2").Value = $newZhHandoffDate

$rZhA2 = $wsZh.Range("A2")
$rZhA2.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($rZhA2, $oldHyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newMdName) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de": A2 file name (hyperlink), G2 latest handoff xlf file,
# H2 latest handoff datetime (shares the same text as Overview!G2)
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("G2").Value = $newDeXlf
$wsDe.Range("H2").Value = $newHandoffDate

$rDeA2 = $wsDe.Range("A2")
$rDeA2.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($rDeA2, $oldHyperlinkAddress, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $newMdName) | Out-Null

Write-Host "Report regenerated for handoff."
